# Daily attendance processing - 2026-01-05 19:36:16
# Normalize the "Recorded By" (column G) values: move the "System" token
# that currently leads the comma-separated list to the position dictated
# by the latest processing run (System no longer sorts first for these
# specific combinations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after replacements observed for the "Recorded By" column.
$replacements = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com"
}

# Determine the last used row from the worksheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Text               # reading via .Value is unreliable on this host; .Text works

    if ($null -ne $value -and $replacements.ContainsKey($value)) {
        $cell.Value = $replacements[$value]
    }
}
